$wb = $excel.ActiveWorkbook

# --- New sheet: COLOR CODING (added after the last existing sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws6.Name = "COLOR CODING"

$ws6.Range("A1").Value = "PLATE NUMBER"
$ws6.Range("B1").Value = "CLASS"
$ws6.Range("C1").Value = "COLOR"
$ws6.Range("D1").Value = "DATE"
$ws6.Range("E1").Value = "TIME"

$ws6.Range("A2").Value = "abc"
$ws6.Range("B2").Value = "truck"
$ws6.Range("C2").Value = "red"
$ws6.Range("D2").Value = "Jan 12, 2015"
$ws6.Range("E2").Value = "12:00NN"

$ws6.Range("A3").Value = "abc"
$ws6.Range("B3").Value = "truck"
$ws6.Range("C3").Value = "red"
$ws6.Range("D3").Value = "Dec 13, 2016"
$ws6.Range("E3").Value = "11:00AM"

$ws6.Columns.Item(1).ColumnWidth = 23.42578125
$ws6.Columns.Item(2).ColumnWidth = 18.7109375
$ws6.Columns.Item(3).ColumnWidth = 16.85546875
$ws6.Columns.Item(4).ColumnWidth = 15.7109375
$ws6.Columns.Item(5).ColumnWidth = 21.140625

$ws6.Range("D7").Select()

# --- SPEEDING (sheet1): add a new row 4 ---
$ws1 = $wb.Worksheets.Item("SPEEDING")
$ws1.Range("A4").Value = "xyz"
$ws1.Range("B4").Value = "truck"
$ws1.Range("C4").Value = "black"
$ws1.Range("D4").Value = "July 3, 2016"
$ws1.Range("E4").Value = "3:00am"
$ws1.Range("A4:E4").Select()

# --- SWERVING (sheet2): update row 2 DATE/TIME ---
$ws2 = $wb.Worksheets.Item("SWERVING")
$ws2.Range("D2").Value = "November 3, 2016"
$ws2.Range("E2").Value = "3:00am"
$ws2.Range("E13").Select()

# --- DRUNK DRIVING (sheet3): no data change, just tab selection moves away ---
$ws3 = $wb.Worksheets.Item("DRUNK DRIVING")

# --- COUNTERFLOWING (sheet4): add new row 2, becomes active tab ---
$ws4 = $wb.Worksheets.Item("COUNTERFLOWING")
$ws4.Range("A2").Value = "xyz"
$ws4.Range("B2").Value = "truck"
$ws4.Range("C2").Value = "black"
$ws4.Range("D2").Value = "September 3, 2016"
$ws4.Range("E2").Value = "3:00am"
$ws4.Range("E7").Select()

# --- BEATING THE RED LIGHT (sheet5): selection change only ---
$ws5 = $wb.Worksheets.Item("BEATING THE RED LIGHT")
$ws5.Range("A1:E1").Select()

$ws4.Activate()
